# Fix import functionality and add newAltIndexToAdd array
#
# The "Alt Index" lookup column on "Table 1" held a handful of rows whose
# generated reference codes (R-218.../R-222...) were stale/out of order
# relative to the new "AltIndex-z" marker row. Re-enter the corrected
# values so the workbook (and its shared-string table) reflect the fixed
# import data, and restore the active selection to the cell the user was
# last working on.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Table 1")

# Re-enter the corrected Alt Index reference codes for rows 19-23.
$ws.Range("A19").Value = "R-218-I18-Cf1-AE28"
$ws.Range("A20").Value = "R-219-I19-Ce1-8FFF"
$ws.Range("A21").Value = "R-220-I20-CT1-6919"
$ws.Range("A22").Value = "R-221-I21-CE2-9AD5"
$ws.Range("A23").Value = "R-222-I22-CS2-06EC"

# Move/restore the active selection to where the user left off.
$null = $ws.Range("H20").Select()
